$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Qty-executed-up-to-date (column C) updates ---
$ws.Range("C8").Value = 67
$ws.Range("C9").Value = 95
$ws.Range("C10").Value = 65
$ws.Range("C11").Value = 58
$ws.Range("C12").Value = 26
$ws.Range("C13").Value = 51
$ws.Range("C14").Value = 58
$ws.Range("C15").Value = 7
$ws.Range("C16").Value = 70
$ws.Range("C17").Value = 70

# --- Upto date Amount (column G) updates - stored as text e.g. "24320.00" ---
# Force text storage (matches existing cells, which are text-typed "X.00" strings)
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "24320.00"

$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "30680.00"

$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "38396.00"

$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "6936.00"

$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "1334.00"

# --- Grand Total rows (19 and 21): G and H columns ---
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "101666.00"

$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "101666.00"

$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "101666.00"

$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "101666.00"
